# Update Price (D) and Volume(1h) (E) columns for the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.435.03"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "'1.678.21"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'216.90"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.2699"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").Value = "'0.06401"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("D11").Value = "'0.07815"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "'4.514"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "'1.673.51"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "'0.5565"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'0.0₅8337"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "'65.62"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'26.484.19"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "'193.44"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'142.18"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").Value = "'0.1288"
$ws.Range("E25").Value = "  +5.64%  "
$ws.Range("D26").Value = "'7.407"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'16.26"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").Value = "'1.437"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "'0.06280"
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("D30").Value = "'1.274"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "'3.613"
$ws.Range("E31").Value = "  +4.64%  "
$ws.Range("D32").Value = "'3.451"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").Value = "'1.680"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "'0.6121"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").Value = "'2.785"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "'0.01633"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'6.137"
$ws.Range("E39").Value = "  +7.28%  "
$ws.Range("D40").Value = "'1.084.96"
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("D41").Value = "'0.8648"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'100.33"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'1.823.27"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").Value = "'57.23"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").Value = "'8.146"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -5.06%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "'1.474"
$ws.Range("E50").Value = "  +5.78%  "
$ws.Range("D51").Value = "'6.027"
$ws.Range("E51").Value = "  +1.76%  "
